# Rename the diff-header columns so they carry the respective AHB
# format-version suffix (FV2310 / FV2404) instead of the generic
# "_old" / "_new" suffixes, then turn the header row into a proper
# Excel Table and freeze it in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "old" (left-hand / FV2310) header block: columns A-J ---
$fv2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $fv2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310[$i]
}

# Column K1 stays "diff" (untouched).

# --- 2. Rename the "new" (right-hand / FV2404) header block: columns L-U ---
$fv2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $fv2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404[$i]
}

# --- 3. Turn the data range into an Excel Table ("Table1") ---
$usedRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"

# --- 4. Freeze the header row (split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
